$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name swap (rows 36 and 37) ---
$ws.Range('B36').Value = 'Hedera'
$ws.Range('B37').Value = 'WEMIXToken'

# --- Link swap (rows 36 and 37) ---
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

# --- Price (column D) updates; force text format so numeric-looking
#     strings (e.g. "304.78") are kept as text, not converted to numbers ---
$dFullRange = $ws.Range("D2:D51")
$dFullRange.NumberFormat = "@"
$ws.Range('D2').Value = '43.177.20'
$ws.Range('D3').Value = '2.328.95'
$ws.Range('D5').Value = '304.78'
$ws.Range('D6').Value = '97.75'
$ws.Range('D11').Value = '19.51'
$ws.Range('D14').Value = '6.95'
$ws.Range('D15').Value = '2.692.22'
$ws.Range('D16').Value = '2.345.20'
$ws.Range('D18').Value = '43.062.79'
$ws.Range('D19').Value = '12.63'
$ws.Range('D22').Value = '67.99'
$ws.Range('D23').Value = '237.91'
$ws.Range('D24').Value = '2.21'
$ws.Range('D27').Value = '24.97'
$ws.Range('D28').Value = '166.12'
$ws.Range('D30').Value = '9.15'
$ws.Range('D31').Value = '33.29'
$ws.Range('D33').Value = '18.08'
$ws.Range('D35').Value = '4.56'
$ws.Range('D36').Value = '0.0698'
$ws.Range('D37').Value = '2.35'
$ws.Range('D39').Value = '2.81'
$ws.Range('D41').Value = '0.109'
$ws.Range('D42').Value = '1.994.25'
$ws.Range('D43').Value = '10.73'
$ws.Range('D45').Value = '18.16'
$ws.Range('D47').Value = '2.79'
$ws.Range('D49').Value = '2.558.22'
$ws.Range('D50').Value = '53.70'
$ws.Range('D51').Value = '72.05'
$dFullRange.Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('E7').Value = '  -1.29%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -0.76%  '
$ws.Range('E10').Value = '  -1.03%  '
$ws.Range('E11').Value = '  +6.84%  '
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('E14').Value = '  +1.89%  '
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('E17').Value = '  +0.87%  '
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('E24').Value = '  +2.70%  '
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('E29').Value = '  +1.97%  '
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('E33').Value = '  +5.71%  '
$ws.Range('E34').Value = '  -0.93%  '
$ws.Range('E35').Value = '  -8.61%  '
$ws.Range('E36').Value = '  +1.43%  '
$ws.Range('E37').Value = '  -1.27%  '
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('E39').Value = '  +2.23%  '
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('E41').Value = '  -0.54%  '
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('E43').Value = '  +5.80%  '
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('E45').Value = '  +3.46%  '
$ws.Range('E46').Value = '  -3.94%  '
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('E48').Value = '  -6.14%  '
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('E51').Value = '  -0.76%  '
